# "Final edit for the database sample"
# Adds a new `current_location_id` column (H) to the `users` sheet, with a
# value per user that references a location from `locations`/`stock_levels`,
# formats the new header cell like the other header cells, widens the
# column, and leaves the `users` sheet active/selected (matching the last
# thing the author touched before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# --- new header cell --------------------------------------------------
$ws.Range("H1").Value = "current_location_id"

# Match the look of the existing header row (bold font, thin box border,
# centered/top aligned) by copying the formatting from an existing header
# cell, then nudge the fill so Excel records a distinct style for it.
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Interior.ColorIndex = -4142

# --- new column data ----------------------------------------------------
$ws.Range("H2").Value = 501
$ws.Range("H3").Value = 501
$ws.Range("H4").Value = 502
$ws.Range("H5").Value = 502
$ws.Range("H6").Value = 503
$ws.Range("H7").Value = 501
$ws.Range("H8").Value = 503
$ws.Range("H9").Value = 502
$ws.Range("H10").Value = 501

# --- column width ---------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 21.7

# --- leave the sheet the way the author left it: `users` active, with
#     H15 selected -----------------------------------------------------
$ws.Activate()
$ws.Range("H15").Select()
